$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 0.66008061520584582
$ws.Range("C1").Value = 0.85773800948435164
$ws.Range("C2").Value = 0.94318308348722735
$ws.Range("AA2").Value = 0.91469048864903979
$ws.Range("C4").Value = 0.95140298824639147
$ws.Range("BC4").Value = 0.87650414060866744
$ws.Range("C5").Value = 0.87219328746551494
$ws.Range("G5").Value = 0.71024427908503629
$ws.Range("AT5").Value = 0.86550479276911685
$ws.Range("D6").Value = 0.97157475064734933
$ws.Range("G6").Value = 0.90251529919158058
$ws.Range("K7").Value = 0.84335431426673746
$ws.Range("AC7").Value = 0.73686632960204901
$ws.Range("F8").Value = 0.82707149834508364
$ws.Range("I8").Value = 0.79075219900630711
$ws.Range("J9").Value = 0.91611015736542178
$ws.Range("K9").Value = 0.84516456812807816
$ws.Range("BA9").Value = 0.91140085600005294
$ws.Range("K10").Value = 0.78108678385852204
$ws.Range("M11").Value = 0.89939500789741533
$ws.Range("J12").Value = 0.83787369616884866
$ws.Range("M12").Value = 0.97079350403816012
$ws.Range("AZ12").Value = 0.95549117600007105
$ws.Range("N13").Value = 0.97010106061156109
$ws.Range("O14").Value = 0.92682069993854199
$ws.Range("P14").Value = 0.67589521676283837
$ws.Range("M15").Value = 0.73816405289401599
$ws.Range("Q15").Value = 0.97640624135201215
$ws.Range("AF15").Value = 0.70030580139753118
$ws.Range("O16").Value = 0.55933824472668581
$ws.Range("Q16").Value = 0.85648324997055603
$ws.Range("R16").Value = 0.98633900993852697
$ws.Range("R17").Value = 0.8124625174102359
$ws.Range("AJ18").Value = 0.6910857498080103
$ws.Range("Q19").Value = 0.99919270723377696
$ws.Range("R19").Value = 0.66002020419363516
$ws.Range("T19").Value = 0.62832311926778828
$ws.Range("U20").Value = 0.65362252486644135
$ws.Range("V20").Value = 0.68934491240258122
$ws.Range("S21").Value = 0.98407997302938488
$ws.Range("V21").Value = 0.77813169114275071
$ws.Range("W21").Value = 0.84493792793640976
$ws.Range("G22").Value = 0.85464342111928659
$ws.Range("V24").Value = 0.9987105092955304
$ws.Range("W24").Value = 0.94949265876018285
$ws.Range("Y24").Value = 0.77691204487638799
$ws.Range("W25").Value = 0.76961044777185894
$ws.Range("AA25").Value = 0.72873224602605935
$ws.Range("X26").Value = 0.76879658354274449
$ws.Range("Y26").Value = 0.70089874741032432
$ws.Range("Z27").Value = 0.66603033326119254
$ws.Range("AB27").Value = 0.87922906064127315
$ws.Range("AC27").Value = 0.8476962575626984
$ws.Range("G28").Value = 0.87809596028241665
$ws.Range("AD28").Value = 0.7332383303780754
$ws.Range("AB29").Value = 0.80808158554570841
$ws.Range("AD29").Value = 0.98257445075793359
$ws.Range("AF30").Value = 0.90503029058049567
$ws.Range("AD31").Value = 0.7324032496741455
$ws.Range("AE32").Value = 0.63406876251720645
$ws.Range("AH32").Value = 0.74622784368888262
$ws.Range("AE33").Value = 0.96114792428486118
$ws.Range("AH33").Value = 0.98129652867029005
$ws.Range("H34").Value = 0.58469688988581359
$ws.Range("AJ34").Value = 0.74619131055318899
$ws.Range("U35").Value = 0.9781504056856094
$ws.Range("AG35").Value = 0.53677065308504246
$ws.Range("L36").Value = 0.91082050302070217
$ws.Range("AI36").Value = 0.77428685546291232
$ws.Range("AT36").Value = 0.9505975303558063
$ws.Range("AL37").Value = 0.55093205151933522
$ws.Range("AM37").Value = 0.74657435672191141
$ws.Range("K38").Value = 0.95064184522342154
$ws.Range("AJ38").Value = 0.91101428146140073
$ws.Range("AL39").Value = 0.93013820724350715
$ws.Range("AN39").Value = 0.79607388646130905
$ws.Range("AL40").Value = 0.87858608227682855
$ws.Range("AO40").Value = 0.85675857172859637
$ws.Range("AM41").Value = 0.93411881220491033
$ws.Range("AP41").Value = 0.78209121309519014
$ws.Range("AN42").Value = 0.65424595217635595
$ws.Range("AQ42").Value = 0.91080804767961354
$ws.Range("AR42").Value = 0.94697486199727932
$ws.Range("AT42").Value = 0.61451209879867053
$ws.Range("AO43").Value = 0.93282056321593465
$ws.Range("AR43").Value = 0.9143195548519345
$ws.Range("AS43").Value = 0.90440335651761761
$ws.Range("AR45").Value = 0.90614260799800617
$ws.Range("AU45").Value = 0.85003344276697135
$ws.Range("AR46").Value = 0.7476341779693032
$ws.Range("AS46").Value = 0.86883069544116731
$ws.Range("AU46").Value = 0.86067322076434616
$ws.Range("AV47").Value = 0.67399776886403129
$ws.Range("AX48").Value = 0.96371323655344932
$ws.Range("AO49").Value = 0.65075319220484962
$ws.Range("AU49").Value = 0.94050609066832946
$ws.Range("AZ50").Value = 0.85990655684778849
$ws.Range("AW51").Value = 0.74941017511008501
$ws.Range("AY52").Value = 0.93841320946741824
$ws.Range("BB52").Value = 0.9920648886834682
$ws.Range("Z53").Value = 0.78811044897191063
$ws.Range("AP53").Value = 0.96855056491321068
$ws.Range("AY53").Value = 0.76098265170202883
$ws.Range("BD54").Value = 0.90960830310596097
$ws.Range("H55").Value = 0.79897684353392817
$ws.Range("AQ55").Value = 0.99454108088488724
$ws.Range("BB55").Value = 0.95988510384253534
$ws.Range("BK55").Value = 0.65276789732149232
$ws.Range("AV56").Value = 0.88261770116508476
$ws.Range("BE56").Value = 0.96615368944914393
$ws.Range("BF56").Value = 0.82135627896165531
$ws.Range("BE58").Value = 0.99606756933614882
$ws.Range("BF59").Value = 0.7533504603865081
$ws.Range("BH59").Value = 0.96965057001233634
$ws.Range("BI59").Value = 0.64251838548036
$ws.Range("AI60").Value = 0.86099041699256029
$ws.Range("BF60").Value = 0.76125823134898973
$ws.Range("BI60").Value = 0.79164086803453182
$ws.Range("U61").Value = 0.995079797535243
$ws.Range("Q62").Value = 0.77905617660878712
$ws.Range("BE62").Value = 0.87216306100981433
$ws.Range("BI62").Value = 0.74032461103432534
$ws.Range("BK62").Value = 0.92216323753350249
$ws.Range("BN64").Value = 0.80380823126328682
$ws.Range("B65").Value = 0.81008797544691169
$ws.Range("AX65").Value = 0.66257813700878643
$ws.Range("BK65").Value = 0.82389120912184644
$ws.Range("BL65").Value = 0.9552637174364853
$ws.Range("Q66").Value = 0.78413301549097381
$ws.Range("BM66").Value = 0.99143935447944154
$ws.Range("BP66").Value = 0.94663351959065278
$ws.Range("A67").Value = 0.77843747638191996
$ws.Range("N67").Value = 0.6216765521275518
$ws.Range("BK67").Value = 0.88973807547536521
$ws.Range("AJ68").Value = 0.77875499489480693
$ws.Range("BA68").Value = 0.85655977941534933
